$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 287; this pushes all existing rows
# 287..359 down to 288..360 and extends the used range to A1:R360.
$ws.Rows.Item(287).Insert()

# Populate the newly inserted row 287 with the new record.
$ws.Range("A287").Value = 5
$ws.Range("B287").Value = "Macroferia Regional de Talca"
$ws.Range("C287").Value = "Maule"
$ws.Range("D287").Value = 44855
$ws.Range("E287").Value = 7
$ws.Range("F287").Value = 100112008
$ws.Range("G287").Value = "Coliflor"
$ws.Range("H287").Value = "Sin especificar"
$ws.Range("I287").Value = "Primera"
$ws.Range("J287").Value = 3000
$ws.Range("K287").Value = 900
$ws.Range("L287").Value = 900
$ws.Range("M287").Value = 900
$ws.Range("N287").Value = "$/unidad"
$ws.Range("O287").Value = "Región del Maule"
$ws.Range("P287").Value = 900
$ws.Range("Q287").Value = 1
$ws.Range("R287").Value = "Hortaliza"
